# Shane Watson vs Chennai Super Kings sheet — the per-innings rows (runs,
# balls, fours, sixes in columns C:F) were reshuffled into a different row
# order. Re-write each row's C:F with its new target figures. The sheet
# stores these numeric-looking values as text (see the numberStoredAsText
# ignoredError on the sheet), so every value is written with a leading
# apostrophe to force a text/quote-prefixed cell instead of a numeric one —
# exactly like the pre-existing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowStats($row, $runs, $balls, $fours, $sixes) {
    $ws.Cells.Item($row, 3).Value = "'" + $runs
    $ws.Cells.Item($row, 4).Value = "'" + $balls
    $ws.Cells.Item($row, 5).Value = "'" + $fours
    $ws.Cells.Item($row, 6).Value = "'" + $sixes
}

Set-RowStats 2  14 16 1 1
Set-RowStats 3  50 40 6 1
Set-RowStats 4  83 53 11 3
Set-RowStats 5  1  6  0 0
Set-RowStats 6  42 38 1 3
Set-RowStats 7  36 28 6 0
Set-RowStats 8  14 19 1 1
Set-RowStats 9  4  5  1 0
Set-RowStats 10 33 21 1 4
Set-RowStats 11 14 18 3 0
Set-RowStats 12 8  3  2 0
